# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de report sheets to reflect a newer report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 21:27:18"
$wsZhCn.Range("H2").Value = "2016-03-24 21:27:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 21:27:23"
$wsDeDe.Range("H2").Value = "2016-03-24 21:27:53"
